$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row before row 31 (shifts old rows 31-38 down to 32-39)
$ws.Rows("31:31").Insert()

# 2. Match formatting on the new row 31 and on rows 16/22 to the style
#    already used on row 15 (the highlighted "indicator species" block).
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A31").PasteSpecial(-4122)

$ws.Range("B15:C15").Copy()
$ws.Range("B16:C16").PasteSpecial(-4122)
$ws.Range("B22:C22").PasteSpecial(-4122)
$ws.Range("B31:C31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Mark row 31 as a custom-formatted row (matches rows 28-30 in the same block).
$ws.Rows("31:31").RowHeight = $ws.Rows("30:30").RowHeight

# 3. Populate the new row 31 with the "Extended data table 1 caption" task,
#    matching the same Goal/Suggested-work-plan text used in row 22
#    ("Finish formatting indicator species table for extended data").
$ws.Range("A31").Value = "Extended data table 1 caption"
$ws.Range("B31").Value = "Complete by Friday, August 18th"
$ws.Range("C31").Value = "Will probably do this before then, while taking a break from writing"

# 4. Update the sheet view (scroll position / selection) to match the
#    author's final cursor position.
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("A16:C16").Select()
